$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the old row 10 ("actual_demand" / 1990 entry),
# pushing the 1990-2019 block (and the trailing blank row) down by one.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with a new "cost_variable_om" parameter entry.
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "dem_elec"
$ws.Range("C10").Value = "cost_variable_om"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0

# Move the active selection the same way the author's session ended up.
$ws.Range("E11").Select() | Out-Null

# Grow the autofilter range by one row to cover the newly added row.
$ws.AutoFilterMode = $false
$ws.Range("A5:L603").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the autofilter range.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$5:`$L`$603"
